# Update countries & provincias Spain
# - Refresh the "Datos actualizados" timestamp (09:05 -> 09:35)
# - Refresh a handful of per-country statistics with newer scraped numbers
# - El Salvador's total cases overtook Lituania and Somalia, so those three
#   countries re-sort (El Salvador, Lituania, Somalia) in the ranked table
# - Montserrat overtook Seychelles (tie-break reorder) in the ranked table
# - Santa Lucia / Nueva Caledonia and Bonaire / Sahara Occidental swap places
#   too (their stats are identical so only the labels move)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp banner
$ws.Range("A1").Value = "Datos actualizados a 21 de Mayo de 2020 a las 09:35"

# Polonia (row 34): recuperados / casos activos refreshed
$ws.Range("D34").Value = 8452
$ws.Range("E34").Value = 10325

# Rows 92-94 re-sorted by total cases: El Salvador, Lituania, Somalia
$ws.Range("A92").Value = "El Salvador"
$ws.Range("B92").Value = 1640
$ws.Range("C92").Value = 69
$ws.Range("D92").Value = 540
$ws.Range("E92").Value = 1068
$ws.Range("F92").Value = 0
$ws.Range("G92").Value = 1
$ws.Range("H92").Value = 32

$ws.Range("A93").Value = "Lituania"
$ws.Range("B93").Value = 1577
$ws.Range("C93").Value = 0
$ws.Range("D93").Value = 1049
$ws.Range("E93").Value = 468
$ws.Range("F93").Value = 0
$ws.Range("G93").Value = 0
$ws.Range("H93").Value = 60

$ws.Range("A94").Value = "Somalia"
$ws.Range("B94").Value = 1573
$ws.Range("C94").Value = 0
$ws.Range("D94").Value = 188
$ws.Range("E94").Value = 1324
$ws.Range("F94").Value = 0
$ws.Range("G94").Value = 0
$ws.Range("H94").Value = 61

# Letonia (row 107): new case counts
$ws.Range("B107").Value = 1025
$ws.Range("C107").Value = 9
$ws.Range("E107").Value = 309
$ws.Range("G107").Value = 1
$ws.Range("H107").Value = 22

# Santa Lucia / Nueva Caledonia swap places (rows 197-198, stats identical)
$ws.Range("A197").Value = "Santa Lucia"
$ws.Range("A198").Value = "Nueva Caledonia"

# Rows 209-210 re-sorted: Montserrat, Seychelles
$ws.Range("A209").Value = "Montserrat"
$ws.Range("B209").Value = 11
$ws.Range("C209").Value = 0
$ws.Range("D209").Value = 10
$ws.Range("E209").Value = 0
$ws.Range("F209").Value = 0
$ws.Range("G209").Value = 0
$ws.Range("H209").Value = 1

$ws.Range("A210").Value = "Seychelles"
$ws.Range("B210").Value = 11
$ws.Range("C210").Value = 0
$ws.Range("D210").Value = 11
$ws.Range("E210").Value = 0
$ws.Range("F210").Value = 0
$ws.Range("G210").Value = 0
$ws.Range("H210").Value = 0

# Bonaire, San Eustaquio y Saba / Sahara Occidental swap places (rows 214-215, stats identical)
$ws.Range("A214").Value = "Bonaire, San Eustaquio y Saba"
$ws.Range("A215").Value = "Sahara Occidental"
